$p = $ppt.ActivePresentation
$s = $p.Slides.Item(19)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$t = $tr.Text
$newText = $t + "`rZoom link: https://cuhk.zoom.us/rec/share/YiSXwJu47YbBeE5K7aqwNRm0SzZ1RCo7oPHwwgY297Icq9GWOfUFxGMHNywxIhUh.ymrIjNPCrpbr8aYm?startTime=1679493977000`rPassword: NECKA2c?"
$tr.Text = $newText
Write-Host "New paragraph count:"
Write-Host $tr.Paragraphs().Count
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    Write-Host "Para $i`: $($tr.Paragraphs($i, 1).Text)"
}
